$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the new columns, matching style of existing header row
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in team record values for every data row (2-33)
$ws.Range("AD2:AD33").Value = 64
$ws.Range("AE2:AE33").Value = 51
$ws.Range("AF2:AF33").Value = 0
